$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test results / timestamps for existing scenario rows
$ws.Range("C2").Value = "Passed"
$ws.Range("H2").Value = "27_04_2020--23_34_44 824"

$ws.Range("H3").Value = "27_04_2020--21_35_53 055"

$ws.Range("H4").Value = "27_04_2020--23_31_17 217"

$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "Failed"
$ws.Range("H5").Value = "27_04_2020--23_35_37 790"

$ws.Range("B6").Value = "Yes"
$ws.Range("C6").Value = "Passed"
$ws.Range("H6").Value = "27_04_2020--23_36_08 311"

# Move the active selection to B4 (was B5)
$ws.Range("B4").Select()
